# "Finished reading chapter 1" — record the actual time it took to read
# Sections 1.1-1.9 (row 3) in the "Actual time length to complete" column (C).
# 15 minutes expressed as a fraction of a day, matching the [h]:mm time format
# already used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

$ws.Range("C3").Value = 15/1440

# Leave the cursor where the author ended up after making the edit.
$ws.Range("C4").Select()
